# Nalco price list update (2025-09-17 06:47:19 UTC)
# A new top entry (Sl.no. 6, basic price 278.95, circular date 17-09-2025) is
# published, so every existing data row shifts down by one and a brand new
# row 2 is created to hold the latest entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 2; this pushes the current rows 2-6 down
# to rows 3-7 (values move, but cell formatting inherited from the insert is
# not reliable, so we fix formatting explicitly below).
$ws.Rows("2:2").Insert()

# Copy the formatting (number formats, alignment, borders, font) from row 3
# -- which now holds the data that used to be in row 2 and therefore still
# carries the correct "data row" styling -- down onto the freshly inserted
# row 2, so the new row looks consistent with the rest of the table.
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new row with the latest circular's data.
$ws.Range("A2").Value = 6
$ws.Range("B2").Value = "ALUMINIUM INGOT"
$ws.Range("C2").Value = "IE07"
$ws.Range("D2").Value = 278.95
$ws.Range("E2").Value = "17-09-2025"
$ws.Range("F2").Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-17-09-2025.pdf"

# The row insert does not relocate the existing hyperlinks together with the
# cells they were attached to, so rebuild the hyperlink list from scratch:
# one hyperlink per row, F2 (new) through F7 (the row that used to be F6).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-17-09-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf")

# Hyperlinks.Add applies Excel's built-in "Hyperlink" cell style (blue,
# underlined), which the source sheet does not use. Restore the original
# plain, centered "data row" formatting on the whole Circular Link column by
# copying it from column E (never touched by the hyperlink styling).
$ws.Range("E2").Copy()
$ws.Range("F2:F7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

Write-Host "Nalco price table updated with 2025-09-17 circular (rows now 1-7)."
